# transport_costs_emissions.xlsx — update model "Year" column values
#   2020 -> 2022   (rows 2-127,  first forecast block)
#   2025 -> 2026   (rows 128-253, second forecast block)
# (2030 / 2040 / 2050 blocks are left untouched)
# Also refresh the active cell selection and turn on a basic page setup,
# matching what the author's Excel session left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("costs_emissions")
$ws.Activate()

# --- Year column (F) bulk updates -----------------------------------
$ws.Range("F2:F127").Value = 2022
$ws.Range("F128:F253").Value = 2026

# --- Page setup -------------------------------------------------------
$ws.PageSetup.PaperSize = 9        # xlPaperA4
$ws.PageSetup.Orientation = 1      # xlPortrait

# --- Selection ----------------------------------------------------------
$ws.Range("L15").Select()
